$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old used range first so rows beyond row 2 are removed entirely
$ws.UsedRange.Clear()

# Set the new header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# Set the new data row
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 8267
$ws.Range("D2").Value = 0.1287877559661865
